# Montecreto.xlsx: re-center the 7-day moving-incidence window on the last
# day instead of the middle day ("finestra incidenza 7gg centrata su ultimo g").
#
# Column B = "nuovi pos." (new daily positives), starting row 2.
# Column C = "somma mobile 7gg." (7-day rolling sum of B).
# Column D = C scaled per 100,000 inhabitants (population = 916).
#
# Previously C(row) summed a window centered on the row (3 days before the
# row through 3 days after it), leaving the first/last 3 rows blank because
# they lack enough neighbours. Now C(row) sums the trailing 7-day window
# ending on the row itself (6 days before through the row), so it can be
# computed all the way through the last day, while the first 6 rows (which
# don't have 6 prior days of data) become blank instead.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$population = 916
$firstDataRow = 2
$lastDataRow = 184
$windowSize = 7

# First full trailing window ends on the row where (row - windowSize + 1)
# reaches the first data row -- i.e. row 8 (rows 2-7 stay/become blank,
# same count of blank rows as before, just moved from the tail to the head).
$firstFullWindowRow = $firstDataRow + $windowSize - 1

for ($r = $firstDataRow; $r -le $lastDataRow; $r++) {
    if ($r -ge $firstFullWindowRow) {
        $windowStart = $r - $windowSize + 1
        $sum = 0
        for ($k = $windowStart; $k -le $r; $k++) {
            $sum = $sum + $ws.Cells.Item($k, 2).Value2
        }
        $ws.Cells.Item($r, 3).Value = $sum
        $ws.Cells.Item($r, 4).Value = $sum * 100000.0 / $population
    } elseif ($ws.Cells.Item($r, 3).Value2 -ne "") {
        # Not enough trailing history yet -- blank out. Only touch cells
        # that actually hold stale data; leave already-blank cells as-is.
        $ws.Cells.Item($r, 3).ClearContents()
        $ws.Cells.Item($r, 4).ClearContents()
    }
}
